# Added email functionality :-)
#
# Adds a new "emailData" worksheet (admin email / admin pass / recipient
# address) at the end of the workbook, wires up mailto: hyperlinks for the
# two e-mail addresses, adds a spare (empty, but hyperlink-styled) cell
# C2 next to the existing hyperlink on the "openBrowser" sheet, and makes
# "openBrowser" the active tab instead of "RegisterNewUser".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "emailData" sheet, appended after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$emailSheet = $wb.Worksheets.Add($null, $lastSheet)
$emailSheet.Name = "emailData"

# Column A: labels, Column B: values (written column-by-column so the
# shared-string table picks up the labels before the values).
$emailSheet.Range("A1").Value = "admin email"
$emailSheet.Range("A2").Value = "admin pass"
$emailSheet.Range("A3").Value = "To addresses"

$emailSheet.Range("B1").Value = "admin@gmail.com"
$emailSheet.Range("B2").Value = "password123"
$emailSheet.Range("B3").Value = "recipient@gmail.com"

# Hyperlink the two e-mail addresses, then re-apply the "Hyperlink" cell
# style so both cells share the workbook's existing hyperlink style.
$emailSheet.Hyperlinks.Add($emailSheet.Range("B1"), "mailto:admin@gmail.com")
$emailSheet.Hyperlinks.Add($emailSheet.Range("B3"), "mailto:recipient@gmail.com")
$emailSheet.Range("B1").Style = "Hyperlink"
$emailSheet.Range("B3").Style = "Hyperlink"

# Leave the cursor where data entry finished.
$emailSheet.Range("B3").Select()

# ---------------------------------------------------------------------
# 2. openBrowser (sheet1): add an extra styled-but-empty cell C2 next to
#    the existing hyperlink cell B2.
# ---------------------------------------------------------------------
$openBrowser = $wb.Worksheets.Item("openBrowser")
$openBrowser.Range("C2").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 3. Make "openBrowser" the active sheet/tab (was "RegisterNewUser").
# ---------------------------------------------------------------------
$openBrowser.Activate()
